# The commit adds one new weekly price record for "Ají" (Americana (o) variety)
# at Macroferia Regional de Talca. In the source data it was inserted as a new
# row right after the existing "Americana (o)" row (row 100), pushing every
# following row (101-165) down by one (now 102-166).
#
# Reproduce that with a real row insert, which shifts the existing rows down
# and lets us fill the now-empty row 101 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 101; rows 101-165 shift down to 102-166.
$ws.Rows(101).Insert()

# Populate the new row 101 with the new data point.
$ws.Range("A101").Value = 5
$ws.Range("B101").Value = "Macroferia Regional de Talca"
$ws.Range("C101").Value = "Maule"
$ws.Range("D101").Value = 44574
$ws.Range("E101").Value = 7
$ws.Range("F101").Value = 100112021
$ws.Range("G101").Value = "Ají"
$ws.Range("H101").Value = "Americana (o)"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 150
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("M101").Value = 15000
$ws.Range("N101").Value = "`$/caja 14 kilos"
$ws.Range("O101").Value = "Región del Maule"
$ws.Range("P101").Value = 1071
$ws.Range("Q101").Value = 14
$ws.Range("R101").Value = "Hortaliza"
